# Qatar Stars League workbook update
# The underlying data rows were re-ordered/corrected: for several match
# entries the row's match data (columns B..AD) was shifted among a small
# group of rows while the row's own index in column A (and the constant
# Div/Date columns A/C/D) stayed put. This script reproduces that by
# reading each row's original B..AD values first (snapshot) and then
# writing them back out to the rotated destination rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is a cycle of 1-based worksheet row numbers.
# Row cycle[i] receives the original B..AD content of cycle[i+1]
# (wrapping around), i.e. content rotates "backwards" through the list.
$cycles = @(
    ,@(18, 19)
    ,@(42, 43)
    ,@(45, 46)
    ,@(62, 63)
    ,@(68, 69)
    ,@(70, 71)
    ,@(83, 84)
    ,@(96, 97)
    ,@(102, 103)
    ,@(118, 119)
    ,@(122, 124, 123, 125, 127)
)

$firstCol = 2   # column B
$lastCol  = 30  # column AD

foreach ($cycle in $cycles) {
    $n = $cycle.Length

    # Snapshot the original values of every row in this cycle before
    # writing anything, since rows may both be a source and destination.
    $snapshot = @{}
    foreach ($r in $cycle) {
        $rowVals = @()
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $rowVals += , ($ws.Cells.Item($r, $c).Value())
        }
        $snapshot[$r] = $rowVals
    }

    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $cycle[$i]
        $srcRow  = $cycle[($i + 1) % $n]
        $vals = $snapshot[$srcRow]
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $ws.Cells.Item($destRow, $c).Value = $vals[$c - $firstCol]
        }
    }
}
